# Addition of Mapping API
# Populate the DataDictionaryTemplate sheet with a source/target mapping table
# and make it the active sheet.

$wb  = $excel.ActiveWorkbook
$ws3 = $wb.Worksheets.Item("DataDictionaryTemplate")

# Header row
$ws3.Range("A1").Value = "source_table"
$ws3.Range("B1").Value = "source_column"
$ws3.Range("C1").Value = "source_data_type"
$ws3.Range("D1").Value = "source_column_description"

# Style the header row: Consolas font, light-green font color
$hdr = $ws3.Range("A1:D1")
$hdr.Font.Color = 0xA2FCA2
$hdr.Font.Name = "Consolas"

# Data rows
$ws3.Range("A2").Value = "customer_orders"
$ws3.Range("B2").Value = "order_id"
$ws3.Range("C2").Value = "int"
$ws3.Range("D2").Value = "Order Management"

$ws3.Range("A3").Value = "customer_orders"
$ws3.Range("B3").Value = "customer_name"
$ws3.Range("C3").Value = "varchar"
$ws3.Range("D3").Value = "Customer Info"

# Column widths to match the populated content
$ws3.Columns.Item(1).ColumnWidth = 14
$ws3.Columns.Item(2).ColumnWidth = 13.33203125
$ws3.Columns.Item(3).ColumnWidth = 15.6640625

# Select C1 and make this the active/displayed sheet
$ws3.Range("C1").Select()
$ws3.Activate()
